# Update the "Förändrad" (column C) date for every existing data row
# (rows 2 through 269) from 45175 (2023-09-06) to 45177 (2023-09-08).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 269; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45177
}

# Row 269 gains an explicit row height (matches the sibling rows above it).
$ws.Rows.Item(269).RowHeight = 15

# Add the new last row (270) for case "A 41429-2023".
$newRow = 270

$ws.Cells.Item($newRow, 1).Value2 = "A 41429-2023"   # Beteckning
$ws.Cells.Item($newRow, 2).Value2 = 45175            # Datum
$ws.Cells.Item($newRow, 3).Value2 = 45177            # Förändrad
$ws.Cells.Item($newRow, 4).Value2 = "VÄRMLANDS LÄN"  # Län
$ws.Cells.Item($newRow, 5).Value2 = "GRUMS"          # Kommun
# Column F (Markägare) intentionally left blank, as in the source data.
$ws.Cells.Item($newRow, 7).Value2 = 0.8              # Area (ha)
$ws.Cells.Item($newRow, 8).Value2 = 0                # Fridlysta
$ws.Cells.Item($newRow, 9).Value2 = 0                # Signalarter
$ws.Cells.Item($newRow, 10).Value2 = 0               # NT
$ws.Cells.Item($newRow, 11).Value2 = 0               # VU
$ws.Cells.Item($newRow, 12).Value2 = 0               # EN
$ws.Cells.Item($newRow, 13).Value2 = 0               # CR
$ws.Cells.Item($newRow, 14).Value2 = 0               # RE
$ws.Cells.Item($newRow, 15).Value2 = 0               # Rödlistade
$ws.Cells.Item($newRow, 16).Value2 = 0               # Hotade
$ws.Cells.Item($newRow, 17).Value2 = 0               # Alla arter

# Apply the same date number format to B270/C270 as the rest of the column.
$ws.Cells.Item($newRow, 2).NumberFormat = $ws.Cells.Item($newRow - 1, 2).NumberFormat
$ws.Cells.Item($newRow, 3).NumberFormat = $ws.Cells.Item($newRow - 1, 3).NumberFormat

# Column R (Artnamn) stays an empty, word-wrapped cell like the rows above it.
$ws.Cells.Item($newRow, 18).WrapText = $true
